# "Test change of self data item template"
# Adds a new "self" data-item column (column J) to Sheet1's data area,
# writing the numbers 1..10 into J2:J11 (one value per existing data row,
# continuing two rows past the prior last row of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 2
$endRow   = 11
$col      = 10  # column J

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Cells.Item($row, $col).Value = $row - $startRow + 1
}
